# Update Name of Algo
# Applies updated KNN imputation result values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.007
$ws.Range("A3").Value = -21.698
$ws.Range("D3").Value = -7.824000000000001
$ws.Range("E6").Value = 16.391
$ws.Range("D12").Value = -7.243
$ws.Range("A14").Value = -21.559
$ws.Range("E19").Value = 16.411
$ws.Range("A21").Value = -20.469
$ws.Range("A23").Value = -20.428
$ws.Range("D24").Value = -7.239999999999999
$ws.Range("E24").Value = 16.923
$ws.Range("A25").Value = -21.405
$ws.Range("B25").Value = 6.576000000000001
$ws.Range("D25").Value = -8.146000000000001
$ws.Range("A26").Value = -21.385
$ws.Range("B27").Value = 6.298
$ws.Range("A29").Value = -21.219
$ws.Range("E30").Value = 16.381
$ws.Range("B31").Value = 6.105999999999999
$ws.Range("E31").Value = 16.329
$ws.Range("E33").Value = 17.229
$ws.Range("B39").Value = 7.717000000000001
$ws.Range("E42").Value = 16.654
$ws.Range("B48").Value = 5.442
$ws.Range("D50").Value = -7.988999999999999
$ws.Range("B51").Value = 6.152
$ws.Range("B52").Value = 5.951000000000001
$ws.Range("A53").Value = -22.042
$ws.Range("D53").Value = -7.77
$ws.Range("B55").Value = 4.697
$ws.Range("E55").Value = 16.354
$ws.Range("B56").Value = 5.003
$ws.Range("A57").Value = -21.303
$ws.Range("B57").Value = 6.431999999999999
$ws.Range("D57").Value = -8.231999999999999
$ws.Range("E58").Value = 16.691
$ws.Range("A59").Value = -22.266
$ws.Range("D61").Value = -7.811
$ws.Range("D63").Value = -7.359999999999999
$ws.Range("E65").Value = 17.084
$ws.Range("A69").Value = -21.531
$ws.Range("D70").Value = -6.705000000000001
$ws.Range("E70").Value = 17.918
$ws.Range("B73").Value = 6.751
$ws.Range("E75").Value = 16.59
$ws.Range("A79").Value = -21.025
$ws.Range("A83").Value = -22
$ws.Range("E83").Value = 16.826
$ws.Range("D86").Value = -8.120000000000001
$ws.Range("E86").Value = 16.4
$ws.Range("B89").Value = 5.856999999999999
$ws.Range("B90").Value = 5.671
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.154000000000001
$ws.Range("A93").Value = -21.401
$ws.Range("E96").Value = 16.645
$ws.Range("E97").Value = 16.796
$ws.Range("D98").Value = -8.062999999999999
$ws.Range("D100").Value = -8.145999999999999
$ws.Range("D102").Value = -8.068000000000001
